$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Devices")

# Update test data values (normal load, cable capacitance etc.)
$ws.Range("D8").Value = 257
$ws.Range("E8").Value = 274
$ws.Range("D9").Value = 287
$ws.Range("E9").Value = 315

# Update sheet view: clear the scrolled top-left cell and move the selection
$ws.Activate()
$ws.Range("E10").Select()
